# Commit swaps the two embedded DrawingML theme colour palettes that ship
# with the deck:
#   - ppt/theme/theme1.xml (the slide master's theme, i.e. the theme that
#     actually paints every slide) goes from the "Integral" palette to the
#     "Office Theme" palette.
#   - ppt/theme/theme2.xml (the notes master's theme) goes from the
#     "Office Theme" palette to the "Integral" palette.
#
# The font scheme and format scheme (fills/lines/effects) are byte-for-byte
# identical between the two theme parts already, so the only substantive
# content delta is the 12 colour-scheme entries (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink). Apply the new palette to the slide master's
# theme through the standard PowerPoint ThemeColorScheme.Colors(n).RGB COM
# surface - this is the theme that is actually rendered behind every slide,
# so it is the part of the swap that has a visible effect on the deck.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

# "Office Theme" palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink),
# expressed as OLE-packed (0x00BBGGRR) RGB integers for the RGB property.
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $cs.Count; $i++) {
    $cs.Colors($i).RGB = $officeThemeColors[$i - 1]
}
